# Auto-generated edit script: update LevePriceNQ/HQ + currentAveragePrice* derived
# profit figures across all 8 sheets, matching the upstream scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3317.6667
$ws.Range("I40").Value = 2975.5
$ws.Range("K40").Value = 2975.5
$ws.Range("M40").Value = -2800.5

$ws.Range("H43").Value = 11400.125
$ws.Range("I43").Value = 7000
$ws.Range("K43").Value = 7000
$ws.Range("M43").Value = -6931

$ws.Range("H55").Value = 1020.7692
$ws.Range("J55").Value = 1343.8334
$ws.Range("L55").Value = 1343.8334
$ws.Range("N55").Value = -1771.8334

$ws.Range("H86").Value = 9589.9
$ws.Range("I86").Value = 9237.5
$ws.Range("J86").Value = 10999.5
$ws.Range("K86").Value = 9237.5
$ws.Range("L86").Value = 10999.5
$ws.Range("M86").Value = -8114.5
$ws.Range("N86").Value = -13245.5

$ws.Range("H89").Value = 9589.9
$ws.Range("I89").Value = 9237.5
$ws.Range("J89").Value = 10999.5
$ws.Range("K89").Value = 46187.5
$ws.Range("L89").Value = 54997.5
$ws.Range("M89").Value = -40571.5
$ws.Range("N89").Value = -66229.5

$ws.Range("H100").Value = 9931.416999999999
$ws.Range("I100").Value = 2598.5
$ws.Range("K100").Value = 2598.5
$ws.Range("M100").Value = -2057.5

$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 80000
$ws.Range("L114").Value = 80000
$ws.Range("N114").Value = -88678

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""

$ws.Range("H138").Value = 5240.57
$ws.Range("I138").Value = 3620.375
$ws.Range("J138").Value = 5752.2104
$ws.Range("K138").Value = 10861.125
$ws.Range("L138").Value = 17256.6312
$ws.Range("M138").Value = -5721.125
$ws.Range("N138").Value = -27536.6312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 517.7778
$ws.Range("I4").Value = 559.2857
$ws.Range("K4").Value = 559.2857
$ws.Range("M4").Value = -443.2857

$ws.Range("H32").Value = 17477.63
$ws.Range("I32").Value = 17155
$ws.Range("K32").Value = 17155
$ws.Range("M32").Value = -16868

$ws.Range("H42").Value = 49250
$ws.Range("I42").Value = 60666.668
$ws.Range("K42").Value = 60666.668
$ws.Range("M42").Value = -60180.668

$ws.Range("H50").Value = 2370.4
$ws.Range("I50").Value = 450
$ws.Range("J50").Value = 3650.6667
$ws.Range("K50").Value = 450
$ws.Range("L50").Value = 3650.6667
$ws.Range("M50").Value = 264
$ws.Range("N50").Value = -5078.6667

$ws.Range("H56").Value = 15000
$ws.Range("J56").Value = 15000
$ws.Range("L56").Value = 15000
$ws.Range("N56").Value = -16484

$ws.Range("H60").Value = 53500
$ws.Range("J60").Value = 57000
$ws.Range("L60").Value = 57000
$ws.Range("N60").Value = -58466

$ws.Range("H61").Value = 8658.833000000001
$ws.Range("I61").Value = 6012.25
$ws.Range("K61").Value = 6012.25
$ws.Range("M61").Value = -5800.25

$ws.Range("H63").Value = 9476.615
$ws.Range("I63").Value = 7299
$ws.Range("K63").Value = 7299
$ws.Range("M63").Value = -6613

$ws.Range("H66").Value = 9476.615
$ws.Range("I66").Value = 7299
$ws.Range("K66").Value = 36495
$ws.Range("M66").Value = -33063

$ws.Range("H74").Value = 436328.12
$ws.Range("I74").Value = 477621.28
$ws.Range("J74").Value = 2750
$ws.Range("K74").Value = 477621.28
$ws.Range("L74").Value = 2750
$ws.Range("M74").Value = -476747.28
$ws.Range("N74").Value = -4498

$ws.Range("H77").Value = 436328.12
$ws.Range("I77").Value = 477621.28
$ws.Range("J77").Value = 2750
$ws.Range("K77").Value = 2388106.4
$ws.Range("L77").Value = 13750
$ws.Range("M77").Value = -2383738.4
$ws.Range("N77").Value = -22486

$ws.Range("H136").Value = 8658.833000000001
$ws.Range("I136").Value = 6012.25
$ws.Range("K136").Value = 18036.75
$ws.Range("M136").Value = -15486.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 38914.82
$ws.Range("I20").Value = 58427.223
$ws.Range("J20").Value = 3792.5
$ws.Range("K20").Value = 58427.223
$ws.Range("L20").Value = 3792.5
$ws.Range("M20").Value = -58180.223
$ws.Range("N20").Value = -4286.5

$ws.Range("H86").Value = 97130.57000000001
$ws.Range("I86").Value = 1927.625
$ws.Range("J86").Value = 401780
$ws.Range("K86").Value = 1927.625
$ws.Range("L86").Value = 401780
$ws.Range("M86").Value = -804.625
$ws.Range("N86").Value = -404026

$ws.Range("H89").Value = 97130.57000000001
$ws.Range("I89").Value = 1927.625
$ws.Range("J89").Value = 401780
$ws.Range("K89").Value = 9638.125
$ws.Range("L89").Value = 2008900
$ws.Range("M89").Value = -4022.125
$ws.Range("N89").Value = -2020132

$ws.Range("H94").Value = 1152.8695
$ws.Range("J94").Value = 1437.75
$ws.Range("L94").Value = 1437.75
$ws.Range("N94").Value = -2339.75

$ws.Range("H134").Value = 4459.346
$ws.Range("I134").Value = 3401.05
$ws.Range("K134").Value = 10203.15
$ws.Range("M134").Value = -7668.150000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14290765
$ws.Range("I31").Value = 30305330
$ws.Range("J31").Value = 7503.676
$ws.Range("K31").Value = 30305330
$ws.Range("L31").Value = 7503.676
$ws.Range("M31").Value = -30305035
$ws.Range("N31").Value = -8093.676

$ws.Range("H34").Value = 14290765
$ws.Range("I34").Value = 30305330
$ws.Range("J34").Value = 7503.676
$ws.Range("K34").Value = 30305330
$ws.Range("L34").Value = 7503.676
$ws.Range("M34").Value = -30305128
$ws.Range("N34").Value = -7907.676

$ws.Range("H99").Value = 4049
$ws.Range("I99").Value = 3635.1428
$ws.Range("J99").Value = 4773.25
$ws.Range("K99").Value = 3635.1428
$ws.Range("L99").Value = 4773.25
$ws.Range("M99").Value = -2137.1428
$ws.Range("N99").Value = -7769.25

$ws.Range("H126").Value = 4049
$ws.Range("I126").Value = 3635.1428
$ws.Range("J126").Value = 4773.25
$ws.Range("K126").Value = 10905.4284
$ws.Range("L126").Value = 14319.75
$ws.Range("M126").Value = -8435.428400000001
$ws.Range("N126").Value = -19259.75

$ws.Range("H132").Value = 19671.092
$ws.Range("I132").Value = 2751.92
$ws.Range("K132").Value = 8255.76
$ws.Range("M132").Value = -5725.76

$ws.Range("H134").Value = 3450.0212
$ws.Range("I134").Value = 3013.6592
$ws.Range("K134").Value = 9040.9776
$ws.Range("M134").Value = -6505.9776

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8136.2
$ws.Range("I2").Value = 68.166664
$ws.Range("K2").Value = 408.999984
$ws.Range("M2").Value = -295.999984

$ws.Range("H4").Value = 2229024.2
$ws.Range("I4").Value = 688159.1
$ws.Range("K4").Value = 2064477.3
$ws.Range("M4").Value = -2064365.3

$ws.Range("H76").Value = 8000
$ws.Range("I76").Value = 8000
$ws.Range("K76").Value = 24000
$ws.Range("M76").Value = -23617

$ws.Range("H79").Value = 8000
$ws.Range("I79").Value = 8000
$ws.Range("K79").Value = 24000
$ws.Range("M79").Value = -22674

$ws.Range("H107").Value = 1101.25
$ws.Range("I107").Value = 877
$ws.Range("J107").Value = 1213.375
$ws.Range("K107").Value = 2631
$ws.Range("L107").Value = 3640.125
$ws.Range("M107").Value = -711
$ws.Range("N107").Value = -7480.125

$ws.Range("H124").Value = 3339441
$ws.Range("I124").Value = 3339441
$ws.Range("K124").Value = 10018323
$ws.Range("M124").Value = -10013413

$ws.Range("H137").Value = 19597.666
$ws.Range("I137").Value = 99999
$ws.Range("K137").Value = 299997
$ws.Range("M137").Value = -294897

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 1497.5
$ws.Range("I36").Value = 1497.5
$ws.Range("K36").Value = 1497.5
$ws.Range("M36").Value = -1012.5

$ws.Range("H126").Value = 3932.0908
$ws.Range("I126").Value = 3045.6667
$ws.Range("K126").Value = 9137.000100000001
$ws.Range("M126").Value = -6667.000100000001

$ws.Range("H132").Value = 4887.5586
$ws.Range("I132").Value = 3743.7307
$ws.Range("J132").Value = 8605
$ws.Range("K132").Value = 11231.1921
$ws.Range("L132").Value = 25815
$ws.Range("M132").Value = -8701.1921
$ws.Range("N132").Value = -30875

$ws.Range("H136").Value = 76537
$ws.Range("J136").Value = 76537
$ws.Range("L136").Value = 229611
$ws.Range("N136").Value = -234711

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4115.3213
$ws.Range("I22").Value = 2045.1111
$ws.Range("J22").Value = 5095.9473
$ws.Range("K22").Value = 2045.1111
$ws.Range("L22").Value = 5095.9473
$ws.Range("M22").Value = -1750.1111
$ws.Range("N22").Value = -5685.9473

$ws.Range("H27").Value = 4115.3213
$ws.Range("I27").Value = 2045.1111
$ws.Range("J27").Value = 5095.9473
$ws.Range("K27").Value = 2045.1111
$ws.Range("L27").Value = 5095.9473
$ws.Range("M27").Value = -1938.1111
$ws.Range("N27").Value = -5309.9473

$ws.Range("H46").Value = 6481.148
$ws.Range("J46").Value = 9529.235000000001
$ws.Range("L46").Value = 9529.235000000001
$ws.Range("N46").Value = -9905.235000000001

$ws.Range("H132").Value = 6037.645
$ws.Range("I132").Value = 4761.4185
$ws.Range("K132").Value = 14284.2555
$ws.Range("M132").Value = -11754.2555

$ws.Range("H136").Value = 5640.8965
$ws.Range("I136").Value = 4365.362
$ws.Range("K136").Value = 13096.086
$ws.Range("M136").Value = -10546.086

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 11399.5

$ws.Range("H25").Value = 12500
$ws.Range("J25").Value = 12500
$ws.Range("L25").Value = 12500
$ws.Range("N25").Value = -13086

$ws.Range("H132").Value = 3983.575
$ws.Range("I132").Value = 3168.3
$ws.Range("J132").Value = 6429.4
$ws.Range("K132").Value = 9504.900000000001
$ws.Range("L132").Value = 19288.2
$ws.Range("M132").Value = -6974.900000000001
$ws.Range("N132").Value = -24348.2
